# Tambhan verify ddl NAP4 cust detail
#
# Updates the "Fee" sheet's Stampduty/Admin fee figures for several
# AssetPriceInclAccessoryAmount bands and moves the active selection /
# active sheet to reflect where the author was last working (Fee!I23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fee")

# Stampduty fee (K/L) bumped from x00000 to x50000 for most price bands,
# plus bespoke Admin/Stampduty fee changes on rows 23, 41 and 47.
$ws.Range("K18").Value = 650000
$ws.Range("L18").Value = 650000

$ws.Range("K19").Value = 750000
$ws.Range("L19").Value = 750000

$ws.Range("K20").Value = 850000
$ws.Range("L20").Value = 850000

$ws.Range("K21").Value = 850000
$ws.Range("L21").Value = 850000

$ws.Range("K22").Value = 850000
$ws.Range("L22").Value = 850000

$ws.Range("H23").Value = 350000
$ws.Range("K23").Value = 340000
$ws.Range("L23").Value = 680000

$ws.Range("K24").Value = 650000
$ws.Range("L24").Value = 650000

$ws.Range("K25").Value = 750000
$ws.Range("L25").Value = 750000

$ws.Range("K26").Value = 850000
$ws.Range("L26").Value = 850000

$ws.Range("K27").Value = 850000
$ws.Range("L27").Value = 850000

$ws.Range("K28").Value = 850000
$ws.Range("L28").Value = 850000

$ws.Range("K29").Value = 850000
$ws.Range("L29").Value = 850000

$ws.Range("K30").Value = 650000
$ws.Range("L30").Value = 650000

$ws.Range("K31").Value = 750000
$ws.Range("L31").Value = 750000

$ws.Range("K32").Value = 850000
$ws.Range("L32").Value = 850000

$ws.Range("K33").Value = 850000
$ws.Range("L33").Value = 850000

$ws.Range("K34").Value = 850000
$ws.Range("L34").Value = 850000

$ws.Range("K35").Value = 850000
$ws.Range("L35").Value = 850000

$ws.Range("K36").Value = 650000
$ws.Range("K37").Value = 750000
$ws.Range("K38").Value = 850000
$ws.Range("K39").Value = 850000
$ws.Range("K40").Value = 850000

$ws.Range("H41").Value = 100000
$ws.Range("I41").Value = 200000
$ws.Range("K41").Value = 300000
$ws.Range("L41").Value = 400000

$ws.Range("K42").Value = 650000
$ws.Range("K43").Value = 750000
$ws.Range("K44").Value = 850000
$ws.Range("K45").Value = 850000
$ws.Range("K46").Value = 850000

$ws.Range("H47").Value = 200000
$ws.Range("I47").Value = 400000
$ws.Range("K47").Value = 600000
$ws.Range("L47").Value = 810000

# Move focus to the Fee sheet (becomes the active tab / sheet) and select
# the cell the author was last looking at.
$ws.Activate()
$ws.Range("I23").Select()
